# Updated symbol list on Tue Dec 13 11:39:26 UTC 2022 with GitHub Actions
#
# Writes numeric-looking values into column D (and a couple of text swaps in
# B/C/E) as literal TEXT, matching the workbook's existing inlineStr storage
# (the sheet stores prices/volumes as text, not numbers). A leading
# apostrophe forces Excel to keep the literal as text instead of coercing it
# to a number; ClearFormats() afterwards drops the transient quote-prefix
# cell format so no stray style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.ClearFormats()
}

# --- Column D price/volume updates -----------------------------------
Set-TextValue "D2"  "268.27"
Set-TextValue "D3"  "21.55"
Set-TextValue "D4"  "6.244"
Set-TextValue "D5"  "0.06175"
Set-TextValue "D6"  "3.566"
Set-TextValue "D7"  "6.554"
Set-TextValue "D8"  "1.365"
Set-TextValue "D9"  "0.8230"
Set-TextValue "D10" "0.01346"
Set-TextValue "D11" "0.1558"
Set-TextValue "D12" "0.08154"
Set-TextValue "D13" "0.03308"
Set-TextValue "D14" "0.03180"
Set-TextValue "D15" "0.09273"
Set-TextValue "D16" "3.746"
Set-TextValue "D17" "0.001626"
Set-TextValue "D18" "0.04673"
Set-TextValue "D19" "0.006395"

Set-TextValue "D23" "3.720"
Set-TextValue "D24" "2.431"
Set-TextValue "D25" "0.3304"

# Row 28 (UpBots) Volume(1h) label gained a "Worstin24h" suffix
Set-TextValue "E28" "27UpBotsUBXTWorstin24h"

Set-TextValue "D40" "0.04661"
Set-TextValue "D41" "0.007000"

# Rows 42/43 swapped places (BKEXToken <-> CEJI) with refreshed data
Set-TextValue "B42" "CEJI"
Set-TextValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003895"
Set-TextValue "E42" "41CEJICEJI"

Set-TextValue "B43" "BKEXToken"
Set-TextValue "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1131"
Set-TextValue "E43" "42BKEXTokenBKK"

Set-TextValue "D44" "0.01187"
Set-TextValue "D45" "0.00006078"
Set-TextValue "D46" "0.0009889"
Set-TextValue "D48" "0.7812"
Set-TextValue "D49" "0.002438"
Set-TextValue "D50" "0.00001898"
Set-TextValue "D51" "0.01239"
